# Edit line graph code - Edit and added correlation for loop codes
#
# The year header row (row 1, columns E:BL) currently stores each year as a
# text label like "1960 [YR1960]". Replace these with plain numeric year
# values (1960-2019), left aligned, so the years can be used numerically
# (e.g. for line-graph / correlation calculations).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCol = 5   # column E
$endCol   = 64  # column BL
$startYear = 1960

for ($col = $startCol; $col -le $endCol; $col++) {
    $year = $startYear + ($col - $startCol)
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $year
    $cell.HorizontalAlignment = -4131  # xlLeft
}

# Reselect the header year range, matching the updated selection left by the edit.
$ws.Range("E1:BL1").Select()
